# Task Breakdown.docx
# Commit: "Added Workspace and base files"
#
# 1. "(Server and Error sim basic UI)" gains " that only implements a quit
#    command" before the closing paren.
# 2. "Shutdown Protocols" gains " (Executed by the UIs)"; the document's
#    lone _GoBack bookmark (previously sitting alone on the last, otherwise
#    empty paragraph) moves to just before that new closing paren.
# 3. A new paragraph "[Asfaqul and Jake]" follows the Shutdown Protocols line.
# 4. "Server sends Ack (Write) or Data (Read)" gains a "Client/" prefix.
# 5. The paragraph that used to hold a single space now reads
#    "[Colin, Jack, and Shaan]".

$d = $word.ActiveDocument

# --- 1: basic UI parenthetical gets the extra clause -----------------------
$rng = $d.Content
$ok1 = $rng.Find.Execute("(Server and Error sim basic UI)", $true, $false, $false, $false, $false, `
                   $true, 1, $false, `
                   "(Server and Error sim basic UI that only implements a quit command)", 2)
Write-Output "edit1 (quit command clause): $ok1"

# --- 2: Shutdown Protocols gets its parenthetical --------------------------
$rng = $d.Content
$ok2 = $rng.Find.Execute("Shutdown Protocols", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Shutdown Protocols (Executed by the UIs)", 2)
Write-Output "edit2 (Executed by the UIs): $ok2"

# --- 3: new paragraph right after it ----------------------------------------
$rng = $d.Content
$ok3 = $rng.Find.Execute("Shutdown Protocols (Executed by the UIs)", $true, $false, $false, $false, $false, `
                   $true, 1, $false, `
                   "Shutdown Protocols (Executed by the UIs)^p[Asfaqul and Jake]", 2)
Write-Output "edit3 ([Asfaqul and Jake] paragraph): $ok3"

# --- move the _GoBack bookmark to sit just before the trailing ")" ---------
# (this also removes it from its old spot, since a document only ever keeps
# one _GoBack bookmark - exactly what the diff shows happening)
$rng = $d.Content
$okF = $rng.Find.Execute("Shutdown Protocols (Executed by the UIs)", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$bmPoint = $d.Range($rng.End - 1, $rng.End - 1)
$d.Bookmarks.Add("_GoBack", $bmPoint)
Write-Output "bookmark relocated: $okF"

# --- 4: "Client/" prefix -----------------------------------------------------
$rng = $d.Content
$ok4 = $rng.Find.Execute("Server sends Ack (Write) or Data (Read)", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Client/Server sends Ack (Write) or Data (Read)", 2)
Write-Output "edit4 (Client/ prefix): $ok4"

# --- 5: turn the lone-space paragraph into the names line -------------------
$ok5 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -eq 2 -and $t.Substring(0, 1) -eq " ") {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $ok5 = $r.Find.Execute(" ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "[Colin, Jack, and Shaan]", 2)
        break
    }
}
Write-Output "edit5 ([Colin, Jack, and Shaan] paragraph): $ok5"
